$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.147202253341675
$ws.Range("B1").Value = 2.443787813186646
$ws.Range("C1").Value = 3.605020523071289
$ws.Range("D1").Value = 1.451034784317017
$ws.Range("E1").Value = 0.9870454668998718
